# Apply odds corrections to Jogos_da_Semana_FlashScore_2024-10-14.xlsx
# Updates specific numeric cells in rows 3, 4, 5, 9 and 10 (betting odds columns)
# on the active worksheet, matching the source data refresh from FlashScore.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 2.1
$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 3
$ws.Range("L3").Value = 4.75
$ws.Range("W3").Value = 5.5
$ws.Range("X3").Value = 8.5
$ws.Range("Z3").Value = 19
$ws.Range("AA3").Value = 23
$ws.Range("AG3").Value = 8
$ws.Range("AI3").Value = 15
$ws.Range("AJ3").Value = 41
$ws.Range("AN3").Value = 4
$ws.Range("AO3").Value = 13
$ws.Range("AW3").Value = 5.5
$ws.Range("AX3").Value = 23

$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("O4").Value = 1.62
$ws.Range("P4").Value = 2.2

$ws.Range("G5").Value = 1.91
$ws.Range("H5").Value = 3.4
$ws.Range("I5").Value = 4.2
$ws.Range("J5").Value = 2.63
$ws.Range("L5").Value = 5
$ws.Range("S5").Value = 1.53
$ws.Range("T5").Value = 2.38
$ws.Range("X5").Value = 8
$ws.Range("Y5").Value = 9.5
$ws.Range("AA5").Value = 19
$ws.Range("AB5").Value = 41
$ws.Range("AC5").Value = 7
$ws.Range("AD5").Value = 6.5
$ws.Range("AG5").Value = 9
$ws.Range("AH5").Value = 19
$ws.Range("AJ5").Value = 41
$ws.Range("AN5").Value = 3.75
$ws.Range("AO5").Value = 11
$ws.Range("AQ5").Value = 41
$ws.Range("AS5").Value = 251
$ws.Range("AT5").Value = 2.38
$ws.Range("AX5").Value = 26

$ws.Range("G9").Value = 2.6
$ws.Range("I9").Value = 2.3
$ws.Range("J9").Value = 3.25
$ws.Range("K9").Value = 2.25
$ws.Range("Q9").Value = 1.73
$ws.Range("R9").Value = 2.08
$ws.Range("S9").Value = 1.33
$ws.Range("T9").Value = 3.25
$ws.Range("U9").Value = 1.62
$ws.Range("V9").Value = 2.2
$ws.Range("W9").Value = 11
$ws.Range("AC9").Value = 13
$ws.Range("AD9").Value = 7
$ws.Range("AL9").Value = 23
$ws.Range("AP9").Value = 21
$ws.Range("AR9").Value = 51
$ws.Range("AS9").Value = 126
$ws.Range("AT9").Value = 3.25
$ws.Range("AV9").Value = 41

$ws.Range("G10").Value = 2.3
$ws.Range("I10").Value = 2.9
$ws.Range("X10").Value = 13
$ws.Range("Y10").Value = 9.5
$ws.Range("AJ10").Value = 29
$ws.Range("AU10").Value = 7
